# Apply the Contoso Chai Tea market trends 2023 edits:
#  1. Rename the "Total Chai sales" header (B1 + the Table1 column header,
#     which is driven by this same cell) to the new Korean label
#     "총 차이 판매액(단위)".
#  2. Change cell D7 from the numeric value 548 to the text value "5:48".
#
# Switch to manual calculation first so existing cached formula results
# (e.g. the shared SUM formula in column B) are left untouched by this
# edit, matching the source diff where only the target cells change.
$excel.Calculation = -4135   # xlCalculationManual

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update the header cell text (this also updates the ListObject/Table
#    column header since the table reads its header names from row 1).
$ws.Range("B1").Value = "총 차이 판매액(단위)"

# 2) Update D7 to hold the text "5:48" instead of the number 548.
$ws.Range("D7").Value = "5:48"
